$wb = $excel.ActiveWorkbook

$newText = "© 2015 - Reviewer Connect`nReviewer Connect Description"

$wsLogin = $wb.Worksheets.Item("LoginHeaderFooter")
$wsLogin.Range("B2").Value = $newText

$wsForgot = $wb.Worksheets.Item("ForgotPasswordHeaderFooter")
$wsForgot.Range("C2").Value = $newText
$wsForgot.Activate()
